$d = $word.ActiveDocument

# 1. Body paragraph: "Date" -> "September 8th, 2024" (with "th" superscripted)
#    Scope the Find to the standalone "Date" paragraph only (there is also a
#    "Date" column header inside the table further down that must stay intact).
$datePara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "Date") {
        $datePara = $p
        break
    }
}

$datePrefix = "September 8"
$supText = "th"
$dateSuffix = ", 2024"
$newDateText = $datePrefix + $supText + $dateSuffix

$datePara.Range.Find.Execute("Date", $true, $false, $false, $false, $false, $true, 1, $false, $newDateText, 2) | Out-Null

# Superscript just the "th" within "September 8th, 2024".
$pStart = $datePara.Range.Start
$supStart = $pStart + $datePrefix.Length
$supEnd = $supStart + $supText.Length
$thRange = $d.Range($supStart, $supEnd)
$thRange.Font.Superscript = $true

# 2. Header: "Student Name 1" -> "Nikhil Chandrashekar"
foreach ($sec in $d.Sections) {
    $hdr = $sec.Headers.Item(1)
    $hdr.Range.Find.Execute("Student Name 1", $true, $false, $false, $false, $false, $true, 1, $false, "Nikhil Chandrashekar", 2) | Out-Null
}
